$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-13 from 45175 to 45183
for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Range("C$row")
    if ($cell.Value2 -eq 45175) {
        $cell.Value2 = 45183
    }
}
